# Revert "update in excel"
# This reverts the earlier edit that:
#   - changed F10 (efficient O) from "O(log(n))" to "O(n)"
#   - changed G10 (highlight) from "Use Binary Search Algorithm along with
#     two pointer technique." to "Use binary search like algorithm with two
#     pointer technique."
#   - added a new row (15) for the "3sum" problem (with its own hyperlink)
# So this script restores F10/G10 to their original text and removes row 15
# (and its hyperlink).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 10 values (efficient O / highlight) ---
$ws.Range("F10").Value = "O(log(n))"
$ws.Range("G10").Value = "Use Binary Search Algorithm along with two pointer technique."

# --- Stash the two "hyperlink" cell formats (Easy/orange = D7, Medium/blue = D3)
#     in scratch cells so we can restore them later after rebuilding the
#     hyperlinks collection (adding a hyperlink resets a cell's style). ---
$ws.Range("D7").Copy() | Out-Null
$ws.Range("ZZ1").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Copy() | Out-Null
$ws.Range("ZZ2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Remove the extra "3sum" row (row 15), shifting rows up ---
$ws.Rows("15:15").Delete()

# --- Rebuild the hyperlinks collection without the removed 3sum link ---
# (the runtime does not clean up hyperlinks automatically when a row is
# deleted, and per-hyperlink .Delete() is a no-op, so delete all and re-add
# the ones that should remain, in original order)
$targets = @(
    @("D2",  "https://leetcode.com/problems/two-sum/description/", "ZZ1"),
    @("D3",  "http://rb.gy/4v8nxh", "ZZ2"),
    @("D4",  "http://rb.gy/oual6", "ZZ2"),
    @("D9",  "http://rb.gy/n391x5", "ZZ2"),
    @("D6",  "http://rb.gy/l4sr2z", "ZZ2"),
    @("D7",  "http://rb.gy/76tm22", "ZZ1"),
    @("D11", "http://rb.gy/1nj72g", "ZZ1"),
    @("D14", "http://rb.gy/nrugfa", "ZZ1"),
    @("D8",  "http://rb.gy/5wfid3", "ZZ1"),
    @("D13", "http://rb.gy/bcqtel", "ZZ1"),
    @("D5",  "http://rb.gy/tqarfs", "ZZ2"),
    @("D12", "https://rebrand.ly/bb0tuzm", "ZZ2"),
    @("D10", "https://rebrand.ly/e3bo9ay", "ZZ2")
)

$ws.Hyperlinks.Delete()
foreach ($pair in $targets) {
    $cell = $pair[0]
    $addr = $pair[1]
    $fmt = $pair[2]
    $ws.Hyperlinks.Add($ws.Range($cell), $addr, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value) | Out-Null
    $ws.Range($fmt).Copy() | Out-Null
    $ws.Range($cell).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# --- Drop the scratch cells used to stash formats ---
$ws.Range("ZZ1").Clear() | Out-Null
$ws.Range("ZZ2").Clear() | Out-Null

# --- Column D width reverts to its original (narrower) size ---
$ws.Columns("D").ColumnWidth = 23.833333333333332

# --- Sheet view: clear the pinned top-left cell and move the selection ---
$ws.Range("D20").Select()
